$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells stay text-typed (matches the original
# inlineStr storage) instead of being auto-coerced to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.214.15"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.38"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.05"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4627"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07268"
$ws.Range("E9").Value = "  -0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8862"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.985.66"
$ws.Range("E11").Value = "  +7.24%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.00"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07797"
$ws.Range("E13").Value = "  -1.79%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.366"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.508"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.14"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008907"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.251.52"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.057"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.048.27"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.949"
$ws.Range("E25").Value = "  +5.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.17"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.45"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.048"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.66"
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.054"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.127"
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7656"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.496"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.749"
$ws.Range("E36").Value = "  +10.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.082"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05239"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01938"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.929"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.063"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5099"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.393"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4782"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.35"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.00"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.636"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06200"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.18"
$ws.Range("E51").Value = "  +0.81%  "
